# Apply the update described by the diff:
# - Row 606: C606 27 -> 29  (B606 formula recalculates to 43476)
# - Row 607: C607 14 -> 21  (B607 formula recalculates to 43497)
# - Row 608: C608 3 -> 66, G608 8 -> 7 (B608 -> 43563, H608 -> 9)
# - Rows 609-611: fill in previously-empty data values so the
#   existing (shared) formulas in B/H/J/K resolve to real numbers
#   instead of the empty string.
#
# Columns L and M in rows 609-611 are formatted as Text (number
# format "@"), so a direct .Value assignment would be stored as a
# text string (matching genuine Excel "typed into a text cell"
# behaviour). The source data instead has real numeric 0s there, so
# we temporarily reset those cells to the default "Normal" style
# before writing the numbers, then restore the original formatting
# (border + text number format) by copying it over from the
# identically-styled cells directly above, which keeps reusing the
# existing style record instead of creating a new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 606
$ws.Range("C606").Value = 29

# Row 607
$ws.Range("C607").Value = 21

# Row 608
$ws.Range("C608").Value = 66
$ws.Range("G608").Value = 7

# Row 609
$ws.Range("C609").Value = 54
$ws.Range("E609").Value = 2
$ws.Range("F609").Value = 2
$ws.Range("G609").Value = 6
$ws.Range("L609:M609").Style = "Normal"
$ws.Range("L609").Value = 0
$ws.Range("M609").Value = 0
$ws.Range("L606:M606").Copy()
$ws.Range("L609:M609").PasteSpecial(-4122)

# Row 610
$ws.Range("C610").Value = 36
$ws.Range("E610").Value = 2
$ws.Range("F610").Value = 2
$ws.Range("G610").Value = 7
$ws.Range("L610:M610").Style = "Normal"
$ws.Range("L610").Value = 0
$ws.Range("M610").Value = 0
$ws.Range("L606:M606").Copy()
$ws.Range("L610:M610").PasteSpecial(-4122)

# Row 611
$ws.Range("C611").Value = 4
$ws.Range("E611").Value = 2
$ws.Range("F611").Value = 2
$ws.Range("G611").Value = 5
$ws.Range("L611:M611").Style = "Normal"
$ws.Range("L611").Value = 0
$ws.Range("M611").Value = 0
$ws.Range("L606:M606").Copy()
$ws.Range("L611:M611").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$excel.CalculateFullRebuild()
$wb.Save()
